$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New job-log rows 34-52 (job numbers 71286-71304), appended below the
# existing data (which ran through row 33).

# A handful of cells hold numeric-looking text (voltages like "460",
# incoming numbers like "82497", unloader counts, and the mm/dd/yyyy
# "enteredOn" dates) that must stay TEXT, matching every other row in
# this sheet (see the numberStoredAsText ignoredError). Mark those cells
# as text *before* writing their value so Excel does not coerce them to
# numbers/dates.
$ws.Range("B40").NumberFormat = "@"
$ws.Range("D38:D41").NumberFormat = "@"
$ws.Range("D43:D49").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E42:E52").NumberFormat = "@"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G45:G50").NumberFormat = "@"
$ws.Range("K34:K52").NumberFormat = "@"

# Row 34 (jobNumber 71286)
$ws.Range("A34").Value = 71286
$ws.Range("B34").Value = '06EY665E 103'
$ws.Range("C34").Value = '0508UE0410'
$ws.Range("D34").Value = 'MULTI'
$ws.Range("E34").Value = '?'
$ws.Range("F34").Value = 'GOOD'
$ws.Range("G34").Value = 'N/A'
$ws.Range("H34").Value = 'NO'
$ws.Range("I34").Formula = '=""'
$ws.Range("J34").Value = 'ravi'
$ws.Range("K34").Value = '9/6/2022'
$ws.Range("L34").Value = $false
$ws.Range("M34").Value = 'N/A'
$ws.Range("N34").Value = 'N/A'
$ws.Range("O34").Value = 'NO'

# Row 35 (jobNumber 71287)
$ws.Range("A35").Value = 71287
$ws.Range("B35").Value = '06E26536A'
$ws.Range("C35").Value = '65334NA2'
$ws.Range("D35").Value = 'MULTI'
$ws.Range("E35").Value = '?'
$ws.Range("F35").Value = 'GOOD'
$ws.Range("G35").Value = 'N/A'
$ws.Range("H35").Value = 'NO'
$ws.Range("I35").Formula = '=""'
$ws.Range("J35").Value = 'ravi'
$ws.Range("K35").Value = '9/6/2022'
$ws.Range("L35").Value = $false
$ws.Range("M35").Value = 'N/A'
$ws.Range("N35").Value = 'N/A'
$ws.Range("O35").Value = 'NO'

# Row 36 (jobNumber 71288)
$ws.Range("A36").Value = 71288
$ws.Range("B36").Value = '06EY77534A'
$ws.Range("C36").Value = '67959NE2'
$ws.Range("D36").Value = 'MULTI'
$ws.Range("E36").Value = '?'
$ws.Range("F36").Value = 'GOOD'
$ws.Range("G36").Value = 'N/A'
$ws.Range("H36").Value = 'NO'
$ws.Range("I36").Formula = '=""'
$ws.Range("J36").Value = 'ravi'
$ws.Range("K36").Value = '9/6/2022'
$ws.Range("L36").Value = $false
$ws.Range("M36").Value = 'N/A'
$ws.Range("N36").Value = 'N/A'
$ws.Range("O36").Value = 'NO'

# Row 37 (jobNumber 71289)
$ws.Range("A37").Value = 71289
$ws.Range("B37").Value = '06ET265360'
$ws.Range("C37").Value = '63522NA0'
$ws.Range("D37").Value = 'MULTI'
$ws.Range("E37").Value = '?'
$ws.Range("F37").Value = 'GOOD'
$ws.Range("G37").Value = 'N/A'
$ws.Range("H37").Value = 'NO'
$ws.Range("I37").Formula = '=""'
$ws.Range("J37").Value = 'ravi'
$ws.Range("K37").Value = '9/6/2022'
$ws.Range("L37").Value = $false
$ws.Range("M37").Value = 'N/A'
$ws.Range("N37").Value = 'N/A'
$ws.Range("O37").Value = 'NO'

# Row 38 (jobNumber 71290)
$ws.Range("A38").Value = 71290
$ws.Range("B38").Value = '06DX3376BC1200'
$ws.Range("C38").Value = '66295NAE2'
$ws.Range("D38").Value = '208'
$ws.Range("E38").Value = '?'
$ws.Range("F38").Value = 'GOOD'
$ws.Range("G38").Value = 'N/A'
$ws.Range("H38").Value = 'NO'
$ws.Range("I38").Formula = '=""'
$ws.Range("J38").Value = 'ravi'
$ws.Range("K38").Value = '9/6/2022'
$ws.Range("L38").Value = $false
$ws.Range("M38").Value = 'N/A'
$ws.Range("N38").Value = 'N/A'
$ws.Range("O38").Value = 'NO'

# Row 39 (jobNumber 71291)
$ws.Range("A39").Value = 71291
$ws.Range("B39").Value = 'DXS45'
$ws.Range("C39").Value = 'XXXXXXX'
$ws.Range("D39").Value = '460'
$ws.Range("E39").Value = '0'
$ws.Range("F39").Value = '?'
$ws.Range("G39").Value = '81605'
$ws.Range("H39").Value = 'NO'
$ws.Range("I39").Value = 'Replaced Input - Brandon'
$ws.Range("J39").Value = 'ravi'
$ws.Range("K39").Value = '9/7/2022'
$ws.Range("L39").Value = $false
$ws.Range("M39").Value = 'N/A'
$ws.Range("N39").Value = 'N/A'
$ws.Range("O39").Value = 'NO'

# Row 40 (jobNumber 71292)
$ws.Range("A40").Value = 71292
$ws.Range("B40").Value = '06E7299610'
$ws.Range("C40").Value = '0305U00659'
$ws.Range("D40").Value = '460'
$ws.Range("E40").Value = '?'
$ws.Range("F40").Value = '?'
$ws.Range("G40").Value = 'N/A'
$ws.Range("H40").Value = 'NO'
$ws.Range("I40").Formula = '=""'
$ws.Range("J40").Value = 'ravi'
$ws.Range("K40").Value = '9/6/2022'
$ws.Range("L40").Value = $false
$ws.Range("M40").Value = 'N/A'
$ws.Range("N40").Value = 'N/A'
$ws.Range("O40").Value = 'NO'

# Row 41 (jobNumber 71293)
$ws.Range("A41").Value = 71293
$ws.Range("B41").Value = '06ET299660'
$ws.Range("C41").Value = '69808ND2'
$ws.Range("D41").Value = '460'
$ws.Range("E41").Value = '?'
$ws.Range("F41").Value = 'GOOD'
$ws.Range("G41").Value = 'N/A'
$ws.Range("H41").Value = 'NO'
$ws.Range("I41").Value = 'Stator only'
$ws.Range("J41").Value = 'ravi'
$ws.Range("K41").Value = '9/6/2022'
$ws.Range("L41").Value = $false
$ws.Range("M41").Value = 'N/A'
$ws.Range("N41").Value = 'N/A'
$ws.Range("O41").Value = 'NO'

# Row 42 (jobNumber 71294)
$ws.Range("A42").Value = 71294
$ws.Range("B42").Value = '6DL3S2700TSK'
$ws.Range("C42").Value = 'NA'
$ws.Range("D42").Value = 'MULTI'
$ws.Range("E42").Value = '0'
$ws.Range("F42").Value = 'GOOD'
$ws.Range("G42").Value = 'N/A82835'
$ws.Range("H42").Value = 'NO'
$ws.Range("I42").Formula = '=""'
$ws.Range("J42").Value = 'ravi'
$ws.Range("K42").Value = '9/7/2022'
$ws.Range("L42").Value = $false
$ws.Range("M42").Value = 'N/A'
$ws.Range("N42").Value = 'N/A'
$ws.Range("O42").Value = 'NO'

# Row 43 (jobNumber 71295)
$ws.Range("A43").Value = 71295
$ws.Range("B43").Value = 'ONE 299 610'
$ws.Range("C43").Value = '408 5J01412'
$ws.Range("D43").Value = '460'
$ws.Range("E43").Value = '0'
$ws.Range("F43").Value = 'GOOD'
$ws.Range("G43").Value = '82585'
$ws.Range("H43").Value = 'NO'
$ws.Range("I43").Value = 'Stator only'
$ws.Range("J43").Value = 'ravi'
$ws.Range("K43").Value = '9/7/2022'
$ws.Range("L43").Value = $false
$ws.Range("M43").Value = 'N/A'
$ws.Range("N43").Value = 'N/A'
$ws.Range("O43").Value = 'NO'

# Row 44 (jobNumber 71296)
$ws.Range("A44").Value = 71296
$ws.Range("B44").Value = 'O6DG5373DC0600'
$ws.Range("C44").Value = '1008U00036'
$ws.Range("D44").Value = '460'
$ws.Range("E44").Value = '2'
$ws.Range("F44").Value = '?'
$ws.Range("G44").Value = 'N/A82920'
$ws.Range("H44").Value = 'NO'
$ws.Range("I44").Value = 'E unl'
$ws.Range("J44").Value = 'ravi'
$ws.Range("K44").Value = '9/7/2022'
$ws.Range("L44").Value = $false
$ws.Range("M44").Value = 'N/A'
$ws.Range("N44").Value = 'N/A'
$ws.Range("O44").Value = 'NO'

# Row 45 (jobNumber 71297)
$ws.Range("A45").Value = 71297
$ws.Range("B45").Value = '06DG5376DC0601'
$ws.Range("C45").Value = '250465 15485'
$ws.Range("D45").Value = '460'
$ws.Range("E45").Value = '2'
$ws.Range("F45").Value = '?'
$ws.Range("G45").Value = '87942'
$ws.Range("H45").Value = 'NO'
$ws.Range("I45").Value = '2  e  unl'
$ws.Range("J45").Value = 'ravi'
$ws.Range("K45").Value = '9/7/2022'
$ws.Range("L45").Value = $false
$ws.Range("M45").Value = 'N/A'
$ws.Range("N45").Value = 'N/A'
$ws.Range("O45").Value = 'NO'

# Row 46 (jobNumber 71298)
$ws.Range("A46").Value = 71298
$ws.Range("B46").Value = '06DG5376DC0601'
$ws.Range("C46").Value = '300019 17066'
$ws.Range("D46").Value = '460'
$ws.Range("E46").Value = '2'
$ws.Range("F46").Value = '?'
$ws.Range("G46").Value = '87941'
$ws.Range("H46").Value = 'NO'
$ws.Range("I46").Value = '2 e unl'
$ws.Range("J46").Value = 'ravi'
$ws.Range("K46").Value = '9/7/2022'
$ws.Range("L46").Value = $false
$ws.Range("M46").Value = 'N/A'
$ws.Range("N46").Value = 'N/A'
$ws.Range("O46").Value = 'NO'

# Row 47 (jobNumber 71299)
$ws.Range("A47").Value = 71299
$ws.Range("B47").Value = 'O6DG5373DC0600'
$ws.Range("C47").Value = '2607U03046'
$ws.Range("D47").Value = '460'
$ws.Range("E47").Value = '2'
$ws.Range("F47").Value = '?'
$ws.Range("G47").Value = '82943'
$ws.Range("H47").Value = 'NO'
$ws.Range("I47").Value = '2 e unl'
$ws.Range("J47").Value = 'ravi'
$ws.Range("K47").Value = '9/7/2022'
$ws.Range("L47").Value = $false
$ws.Range("M47").Value = 'N/A'
$ws.Range("N47").Value = 'N/A'
$ws.Range("O47").Value = 'NO'

# Row 48 (jobNumber 71300)
$ws.Range("A48").Value = 71300
$ws.Range("B48").Value = '06DG53760601'
$ws.Range("C48").Value = 'NA'
$ws.Range("D48").Value = '460'
$ws.Range("E48").Value = '2'
$ws.Range("F48").Value = '?'
$ws.Range("G48").Value = '82944'
$ws.Range("H48").Value = 'NO'
$ws.Range("I48").Value = '2 e unl'
$ws.Range("J48").Value = 'ravi'
$ws.Range("K48").Value = '9/7/2022'
$ws.Range("L48").Value = $false
$ws.Range("M48").Value = 'N/A'
$ws.Range("N48").Value = 'N/A'
$ws.Range("O48").Value = 'NO'

# Row 49 (jobNumber 71301)
$ws.Range("A49").Value = 71301
$ws.Range("B49").Value = '06DG5376DC0601'
$ws.Range("C49").Value = '4619UD4187'
$ws.Range("D49").Value = '460'
$ws.Range("E49").Value = '2'
$ws.Range("F49").Value = '?'
$ws.Range("G49").Value = '82919'
$ws.Range("H49").Value = 'NO'
$ws.Range("I49").Value = '2 e unl'
$ws.Range("J49").Value = 'ravi'
$ws.Range("K49").Value = '9/7/2022'
$ws.Range("L49").Value = $false
$ws.Range("M49").Value = 'N/A'
$ws.Range("N49").Value = 'N/A'
$ws.Range("O49").Value = 'NO'

# Row 50 (jobNumber 71302)
$ws.Range("A50").Value = 71302
$ws.Range("B50").Value = 'O6EY675E 103'
$ws.Range("C50").Value = '2605UE4569'
$ws.Range("D50").Value = 'MULTI'
$ws.Range("E50").Value = '0'
$ws.Range("F50").Value = 'BAD'
$ws.Range("G50").Value = '82251'
$ws.Range("H50").Value = 'NO'
$ws.Range("I50").Formula = '=""'
$ws.Range("J50").Value = 'ravi'
$ws.Range("K50").Value = '9/7/2022'
$ws.Range("L50").Value = $false
$ws.Range("M50").Value = 'N/A'
$ws.Range("N50").Value = 'N/A'
$ws.Range("O50").Value = 'NO'

# Row 51 (jobNumber 71303)
$ws.Range("A51").Value = 71303
$ws.Range("B51").Value = '06CC675E2OO'
$ws.Range("C51").Value = '1701J01641'
$ws.Range("D51").Value = 'MULTI'
$ws.Range("E51").Value = '0'
$ws.Range("F51").Value = 'GOOD'
$ws.Range("G51").Value = 'N/A'
$ws.Range("H51").Value = 'NO'
$ws.Range("I51").Formula = '=""'
$ws.Range("J51").Value = 'ravi'
$ws.Range("K51").Value = '9/7/2022'
$ws.Range("L51").Value = $false
$ws.Range("M51").Value = 'N/A'
$ws.Range("N51").Value = 'N/A'
$ws.Range("O51").Value = 'NO'

# Row 52 (jobNumber 71304)
$ws.Range("A52").Value = 71304
$ws.Range("B52").Value = '6DL32700TSK'
$ws.Range("C52").Value = 'NA'
$ws.Range("D52").Value = 'MULTI'
$ws.Range("E52").Value = '0'
$ws.Range("F52").Value = 'GOOD'
$ws.Range("G52").Value = 'N/A'
$ws.Range("H52").Value = 'NO'
$ws.Range("I52").Formula = '=""'
$ws.Range("J52").Value = 'ravi'
$ws.Range("K52").Value = '9/7/2022'
$ws.Range("L52").Value = $false
$ws.Range("M52").Value = 'N/A'
$ws.Range("N52").Value = 'N/A'
$ws.Range("O52").Value = 'NO'
